# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 17 de Agosto de 2020 a las 10:07"

# --- Update country statistics rows ---

# Row 6: India
$ws.Range("B6").Value = 2651290
$ws.Range("C6").Value = 3974
$ws.Range("D6").Value = 1920265
$ws.Range("E6").Value = 679946
$ws.Range("G6").Value = 34
$ws.Range("H6").Value = 51079

# Row 7: Rusia
$ws.Range("B7").Value = 927745
$ws.Range("C7").Value = 4892
$ws.Range("D7").Value = 736101
$ws.Range("E7").Value = 175904
$ws.Range("G7").Value = 55
$ws.Range("H7").Value = 15740

# Row 48: Singapur
$ws.Range("B48").Value = 55838
$ws.Range("C48").Value = 91
$ws.Range("E48").Value = 3858

# Row 53: Barein
$ws.Range("E53").Value = 3535
$ws.Range("G53").Value = 2
$ws.Range("H53").Value = 172

# Row 57: Armenia
$ws.Range("B57").Value = 41701
$ws.Range("C57").Value = 38
$ws.Range("D57").Value = 34655
$ws.Range("E57").Value = 6222
$ws.Range("G57").Value = 6
$ws.Range("H57").Value = 824

# Row 131: Estonia
$ws.Range("B131").Value = 2192
$ws.Range("C131").Value = 2
$ws.Range("E131").Value = 153

# --- Swap "Islas Malvinas" (row 213) and "Montserrat" (row 214) ---
# Both rows share identical B/C/E/F/G, only the country label, D (Casos activos)
# and H (Muertes) actually move between the two rows.
$ws.Range("A213").Value = "Montserrat"
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0
